# order eims rows & add ncp student creators
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")
$ws.Activate()

# Append new creator rows for NCP student contributors.
$ws.Range("A9").Value = "Arshia"
$ws.Range("C9").Value = "Mehta"
$ws.Range("D9").Value = "Northeast U.S. Shelf LTER"
$ws.Range("E9").Value = "amehta3@wellesley.edu"
$ws.Range("G9").Value = "creator"
$ws.Range("H9").Value = "Northeast U.S. Shelf LTER"
$ws.Range("I9").Value = "NSF"
$ws.Range("J9").Value = "OCE-1655686"

$ws.Range("C10").Value = "Aldrett"
$ws.Range("A10").Value = "Danielle"
$ws.Range("D10").Value = "Northeast U.S. Shelf LTER"
$ws.Range("G10").Value = "creator"
$ws.Range("H10").Value = "Northeast U.S. Shelf LTER"
$ws.Range("I10").Value = "NSF"
$ws.Range("J10").Value = "OCE-1655686"

$ws.Range("J10").Select()
